# Apply "Tried more model optimisation" edit to Model_Optimisation workbook.
# New rows document further ANN / Random forest optimisation iterations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text/content edits, in the order that reproduces the author's shared
#     string table layout (Changes column filled first for rows 6-9, then
#     the Model cell for row 6 switched from "Deep learning" to "ANN",
#     then the remaining Changes column cells for rows 10-14). ---

$ws.Range("C6").Value = "New model type; kept columns from previous iteration"
$ws.Range("C7").Value = "Change data types from float to int to reduce the number of unique values"
$ws.Range("C8").Value = "Drop 'Sex' and 'Stress Level' columns to reduce user input to 8 features"
$ws.Range("C9").Value = "Bin 'age' values "
$ws.Range("B6").Value = "ANN"
$ws.Range("C10").Value = "Bin all columns except for exercise and sedentary hours per week"
$ws.Range("C11").Value = "Same as above but with different model"
$ws.Range("C12").Value = "Reduce n(epochs) to 50 and add a third hidden layer"
$ws.Range("C13").Value = "Increase n(epochs) to 200 and add more neurons to each hidden layer"
$ws.Range("C14").Value = "Use Kerastuner to find the best hyperparameters"

# --- Remaining cells for the new rows (A, B, D columns) ---

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "ANN"
$ws.Range("D7").Value = 0.62

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Random forest"
$ws.Range("D8").Value = 0.62

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "ANN"
$ws.Range("D9").Value = 0.62

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "ANN"
$ws.Range("D10").Value = 0.62

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Random forest"
$ws.Range("D11").Value = 0.63

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "ANN"
$ws.Range("D12").Value = 0.63

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "ANN"
$ws.Range("D13").Value = 0.63

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "ANN"
$ws.Range("D14").Value = 0.653

# Apply number formats: rows 7-13 use the same percent format as existing
# rows (0%), row 14 uses a two-decimal percent format (0.00%)
$ws.Range("D7:D13").NumberFormat = "0%"
$ws.Range("D14").NumberFormat = "0.00%"

# Update selection / active cell to mirror the saved state after edits
$ws.Range("D15").Select()
